$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Capture the existing "premium highlight" (yellow fill) formatting
#    from E2 (the only "Yes" row in the original data) BEFORE we touch
#    any cell content. Stash it on a scratch cell far outside the used
#    range, because once E2's own value/style is overwritten later the
#    original clipboard copy is no longer valid to paste from.
# ---------------------------------------------------------------------
$ws.Range("E2").Copy() | Out-Null
$ws.Range("Z100").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2) Column width changes
#    (COM ColumnWidth reads ~0.83 narrower than the stored OOXML <col> width)
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 54.17   # C: 57 -> 55
$ws.Columns.Item(4).ColumnWidth = 69.17   # D: 60 -> 70
$ws.Columns.Item(8).ColumnWidth = 54.17   # H: 60 -> 55

# ---------------------------------------------------------------------
# 3) Helper to write one data row (columns A-H). Every cell's style is
#    reset to the plain "Normal" style first so stale formatting (e.g.
#    E2's old yellow "Yes" highlight) never survives a value change.
#    NOTE: this interpreter only binds function parameters positionally,
#    so $sheet is passed in explicitly and every call below uses
#    positional (not named) arguments.
# ---------------------------------------------------------------------
function Set-OppRow($sheet, $row, $id, $title, $country, $premium, $applicants, $duration, $org) {
    $sheet.Range("A$row" + ":H$row").Style = "Normal"

    $sheet.Range("A$row").NumberFormat = "@"
    $sheet.Range("A$row").Value = $id
    $sheet.Range("A$row").Style = "Normal"

    $sheet.Range("B$row").Value = "https://aiesec.org/opportunity/global-talent/$id"
    $sheet.Range("C$row").Value = $title
    $sheet.Range("D$row").Value = $country
    $sheet.Range("E$row").Value = $premium
    $sheet.Range("F$row").Value = $applicants
    $sheet.Range("G$row").Value = $duration
    $sheet.Range("H$row").Value = $org
}

# ---------------------------------------------------------------------
# 4) Row data (latest scrape results)
# ---------------------------------------------------------------------
Set-OppRow $ws 2 "1328615" "Back Office Planner" "Madrid, Spain" "No" "10 applicants" "6 - 18 Months" "Mitsubishi Power Europe Sucursal en España"

Set-OppRow $ws 3 "1328614" "Field Service Engineer" "Madrid, Spain" "No" "5 applicants" "6 - 18 Months" "Mitsubishi Power Europe Sucursal en España"

Set-OppRow $ws 4 "1328612" "[DSC] Finance Data Analyst" "Fritz-Erler-Straße 5, 53113 Bonn, Germany" "Yes" "11 applicants" "6 - 18 Months" "DHL Group"

Set-OppRow $ws 5 "1328610" "[CC] Employee Share Plan Support" "Charles-de-Gaulle-Straße 20, 53113 Bonn, Germany" "Yes" "11 applicants" "6 - 18 Months" "DHL Group"

Set-OppRow $ws 6 "1328609" "Taste Hungary | Field Service Support Representative" "Budapeste, Hungria" "Yes" "9 applicants" "9 - 12 Weeks" "EATON"

Set-OppRow $ws 7 "1328556" ": AI Engineer / Machine Learning Developer it" "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt" "No" "7 applicants" "9 - 12 Weeks" "Techno square"

Set-OppRow $ws 8 "1328553" "Business developer" "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt" "No" "1 applicant" "9 - 12 Weeks" "I.C.Robotics"

Set-OppRow $ws 9 "1328551" "Web development" "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt" "No" "2 applicants" "9 - 12 Weeks" "TAR - Company"

Set-OppRow $ws 10 "1328550" "Business development" "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt" "No" "2 applicants" "9 - 12 Weeks" "TAR - Company"

Set-OppRow $ws 11 "1326503" "Data Analyst" "Αθήνα, Ελλάδα" "No" "117 applicants" "9 - 12 Weeks" "Inzeb"

Set-OppRow $ws 12 "1325153" "ACE Program | Danish Language Coach" "Mumbai, Maharashtra, India" "Yes" "3 applicants" "9 - 12 Weeks" "Tata Consultancy Services Ltd."

Set-OppRow $ws 13 "1322596" "HR Intern" "Hyderabad, Telangana, India" "No" "7 applicants" "9 - 12 Weeks" "TERICSOFT TECHNOLOGY SOLUTIONS PVT. LTD."

Set-OppRow $ws 14 "1305239" "Sales" "Eskişehir, Türkiye" "No" "66 applicants" "6 - 18 Months" "Esaysan Endüstriyel Metal Ürünleri Sanayi Ve Ticaret"

# ---------------------------------------------------------------------
# 5) Re-apply the yellow "Yes" premium highlight (captured from the
#    original E2 in step 1, stashed at Z100) to the rows that are "Yes"
#    in the new data, then discard the scratch cell.
# ---------------------------------------------------------------------
$ws.Range("Z100").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null

$ws.Range("Z100").Clear() | Out-Null
